$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DecisionTable")

# C19: "next c" -> "next cc"
$ws.Range("C19").Value = "next cc"

# C20, C21, C22: clear values (set to empty)
$ws.Range("C20").Value = ""
$ws.Range("C21").Value = ""
$ws.Range("C22").Value = ""

# C23: "40" -> "Deleted one condition left"
$ws.Range("C23").Value = "Deleted one condition left"

# C24: "45" -> "Deleted one row below"
$ws.Range("C24").Value = "Deleted one row below"

# Row 25: A25 and C25 become "New Row and Column"; B25, D25, E25 become empty
$ws.Range("A25").Value = "New Row and Column"
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = "New Row and Column"
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = ""
